# Weekly fruit/vegetable price update: insert a new daily record as row 383
# (pushing the existing rows 383-412 down to 384-413) on the single data
# sheet of the "Hortaliza, Vega Monumental Concepción - Brócoli" workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 383..412 down to 384..413, leaving row 383 blank for the new record.
$ws.Rows.Item(383).Insert()

# Populate the newly inserted row 383 with the new weekly price record.
$ws.Cells.Item(383, 1).Value = 11
$ws.Cells.Item(383, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(383, 3).Value = "Bíobío"
$ws.Cells.Item(383, 4).Value = 44931
$ws.Cells.Item(383, 5).Value = 8
$ws.Cells.Item(383, 6).Value = 100112023
$ws.Cells.Item(383, 7).Value = "Brócoli"
$ws.Cells.Item(383, 8).Value = "Sin especificar"
$ws.Cells.Item(383, 9).Value = "Primera"
$ws.Cells.Item(383, 10).Value = 1550
$ws.Cells.Item(383, 11).Value = 750
$ws.Cells.Item(383, 12).Value = 800
$ws.Cells.Item(383, 13).Value = 773
$ws.Cells.Item(383, 14).Value = "$/unidad"
$ws.Cells.Item(383, 15).Value = "Región Metropolitana"
$ws.Cells.Item(383, 16).Value = 773
$ws.Cells.Item(383, 17).Value = 1
$ws.Cells.Item(383, 18).Value = "Hortaliza"
